$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Populate new cells in the same order the original author typed them,
# so new shared-string entries are appended in matching order.
$ws.Range("B11").Value = "E:\Uipath\Emergenteck\My Resumes\0NEW Resume\Vajrang UiPath\Mukesh Kala\HR_GENERATE_OFFER_LETTER_FILES\Offer Letters"
$ws.Range("A11").Value = "OfferLetterFolderPath"
$ws.Range("C11").Value = "Path to store the offer letters Pdfs"

$ws.Range("A13").Value = "EmailSuject"
$ws.Range("C13").Value = "Subject for the e-mail "

$ws.Range("A14").Value = "EmailBody"

$ws.Range("B13").Value = """Congratulations for the offer !"""
$ws.Range("B14").Value = """Congratulations for the offer, Please accept this offer within 10 days and following are terms and conditions ."""

$ws.Range("C14").Value = "E-mail body "

$ws.Range("A16").Value = "CredGmail"
$ws.Range("B16").Value = "GmailSMTP_Credentials"
$ws.Range("C16").Value = "Credential for Gmail "

# Update the selection on the Settings sheet to C6
$ws.Range("C6").Select()
